# Commit: "Add files via upload"
#
# - Adds a new worksheet "x1" right after the existing "default" sheet,
#   fills it with a 12-row English/Thai vocabulary table (A: English,
#   B: Thai), and leaves it as the active/visible tab.
# - On the "default" sheet, moves the selection to A1:B25 and sets the
#   print orientation to portrait.
# - Best-effort: relabel the built-in "Normal" cell style to its Thai
#   localized name (harmless no-op if the host does not expose this).

$wb = $excel.ActiveWorkbook

# --- "default" sheet (existing sheet, stays first / sheet1) ---------------
$ws1 = $wb.Worksheets.Item(1)

# Selection on the default sheet moves to A1:B25. Because the new "x1"
# sheet (added below) becomes the active tab, "default" no longer shows
# as the tab selected when the workbook is reopened.
[void]$ws1.Range("A1:B25").Select()

# Printed page orientation for the default sheet.
$ws1.PageSetup.Orientation = 1

try {
    $wb.Styles.Item(1).Name = "ปกติ"
} catch {
}

# --- new "x1" sheet, inserted right after "default" ------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "x1"

$pairs = @(
    @("enhance", "เพิ่มความสามารถ"),
    @("initiative", "ความคิดริเริ่ม"),
    @("innovate", "นวัตกรรม"),
    @("allocate", "จัดสรร"),
    @("compensate", "ชดเชย"),
    @("coordination", "การประสานงาน"),
    @("implement", "ดำเนินการ"),
    @("incentive", "สิ่งบันเทิง"),
    @("integration", "การรวมกัน"),
    @("legitimate", "ถูกต้อง"),
    @("maximize", "ทำให้เต็มที่"),
    @("objective", "วัตถุประสงค์")
)

for ($i = 0; $i -lt $pairs.Count; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $pairs[$i][0]
    $ws2.Cells.Item($row, 2).Value = $pairs[$i][1]
}

# Selection on the new sheet; it becomes the active/visible tab.
[void]$ws2.Range("K7").Select()
